# Auto-generated by analysis; applies proofErr/run-split restructuring per diff
$d = $word.ActiveDocument
$p2xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="77A100E9" w14:textId="7A3106D0" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:r><w:t>Users (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00732D5B"><w:rPr><w:u w:val="single"/></w:rPr><w:t>user_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">,  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_fname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>use</w:t></w:r><w:r w:rsidR="00DA6745"><w:t>r_l</w:t></w:r><w:r><w:t>name,user_dob</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_paswd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_addr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(2).Range.InsertXML($p2xml)
$p5xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="78A0F08C" w14:textId="66E15150" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:r><w:t xml:space="preserve">Profiles </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">( </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00732D5B"><w:rPr><w:u w:val="single"/></w:rPr><w:t>user</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r w:rsidR="00D17FFE"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>,</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>profile_desciption,profile_photo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(5).Range.InsertXML($p5xml)
$p6xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="5139352D" w14:textId="21363A54" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:r><w:t>Photos (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00732D5B"><w:rPr><w:u w:val="single"/></w:rPr><w:t>photo_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ,</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>photo</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_nme</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>photo_image,photo_datetime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(6).Range.InsertXML($p6xml)
$p7xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="1AB9E52E" w14:textId="1916F809" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:r><w:t>Posts (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00732D5B"><w:rPr><w:u w:val="single"/></w:rPr><w:t>posts_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00D17FFE"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>post_text,datetime</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(7).Range.InsertXML($p7xml)
$p8xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="6899C36A" w14:textId="77777777" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:r><w:t>Groups (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00732D5B"><w:rPr><w:u w:val="single"/></w:rPr><w:t>group_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>id,group</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>_nme,group_description</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) note: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> would be the creator’s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(8).Range.InsertXML($p8xml)
$p9xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="6DEAB211" w14:textId="77777777" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:proofErr w:type="spellStart"/><w:r><w:t>GroupMembers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00732D5B"><w:rPr><w:u w:val="single"/></w:rPr><w:t>group_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>id,user</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>_id</w:t></w:r><w:r><w:t>,member_status</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) note: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>member_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> either being creator, content creator or member</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(9).Range.InsertXML($p9xml)
$p10xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="5FBE19FB" w14:textId="77777777" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:proofErr w:type="spellStart"/><w:r><w:t>GroupPosts</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="002D6D3D"><w:rPr><w:u w:val="single"/></w:rPr><w:t>group_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>id,post</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(10).Range.InsertXML($p10xml)
$p11xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="1AF789ED" w14:textId="38685FD0" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:r><w:t xml:space="preserve">Friends </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">( </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00732D5B"><w:rPr><w:u w:val="single"/></w:rPr><w:t>user</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>_id,friend_id,</w:t></w:r><w:r><w:t>friend_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(11).Range.InsertXML($p11xml)
$p12xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="25731779" w14:textId="5A74E620" w:rsidR="00732D5B" w:rsidRDefault="00732D5B" w:rsidP="00732D5B"><w:proofErr w:type="gramStart"/><w:r><w:t>Comment(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00B8200F"><w:rPr><w:u w:val="single"/></w:rPr><w:t>comment_id,</w:t></w:r><w:r><w:t>post_</w:t></w:r><w:r><w:t>id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00A24B54"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>

'@
$d.Paragraphs(12).Range.InsertXML($p12xml)
Write-Output "done"
